# Insert a new weekly price record as row 40 on the single data sheet.
# This pushes the existing rows 40..105 down to 41..106 (dimension grows
# from A1:R105 to A1:R106) and populates the newly inserted row with the
# new "Haba" record for Terminal Hortofruticola Agro Chillan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40:105 down by inserting a fresh row at position 40.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record.
$ws.Range("A40").Value = 7
$ws.Range("B40").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C40").Value = "Ñuble"
$ws.Range("D40").Value = 45272
$ws.Range("E40").Value = 16
$ws.Range("F40").Value = 100112026
$ws.Range("G40").Value = "Haba"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 12000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = 12000
$ws.Range("N40").Value = "`$/saco 25 kilos"
$ws.Range("O40").Value = "Provincia de Diguillín"
$ws.Range("P40").Value = 480
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
